$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly data rows (2-8) are being re-ordered/rotated: the most recent
# week's data (formerly rows 6-8, date 2021-04-23 / serial 44309) moves to
# the top (rows 2-4), and the older weeks shift down by three rows.
# New row 2-8 values = Old row [6, 7, 8, 2, 3, 4, 5] (in that order).

$ws.Range("D2").Value = 44309
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 305000
$ws.Range("O2").Value = 310000
$ws.Range("P2").Value = 307500
$ws.Range("R2").Value = "Provincia de Cachapoal"
$ws.Range("S2").Value = 683

$ws.Range("D3").Value = 44309
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 285000
$ws.Range("O3").Value = 290000
$ws.Range("P3").Value = 287500
$ws.Range("R3").Value = "Provincia de Cachapoal"
$ws.Range("S3").Value = 639

$ws.Range("D4").Value = 44309
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = 255000
$ws.Range("O4").Value = 260000
$ws.Range("P4").Value = 257500
$ws.Range("R4").Value = "Provincia de Cachapoal"
$ws.Range("S4").Value = 572

$ws.Range("D5").Value = 44295
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 20
$ws.Range("N5").Value = 225000
$ws.Range("O5").Value = 230000
$ws.Range("P5").Value = 227500
$ws.Range("S5").Value = 506

$ws.Range("D6").Value = 44295
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 16
$ws.Range("N6").Value = 195000
$ws.Range("O6").Value = 200000
$ws.Range("P6").Value = 197500
$ws.Range("R6").Value = "Región Metropolitana"
$ws.Range("S6").Value = 439

$ws.Range("D7").Value = 44294
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 20
$ws.Range("N7").Value = 225000
$ws.Range("O7").Value = 230000
$ws.Range("P7").Value = 227500
$ws.Range("R7").Value = "Región Metropolitana"
$ws.Range("S7").Value = 506

$ws.Range("D8").Value = 44294
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 16
$ws.Range("N8").Value = 195000
$ws.Range("O8").Value = 200000
$ws.Range("P8").Value = 197500
$ws.Range("R8").Value = "Región Metropolitana"
$ws.Range("S8").Value = 439
